# Added check for if case statute or degree is none.
# Clear the Degree (column G) for row 6 and the Code/Statute (column F) for row 8,
# since those charge rows have no statute/degree recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("G8").Select()
